# Rename task-order sheets and update their CSV filename references
# to reflect the re-run / re-logged stimulus timing data.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (positional, matches original sheet order) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555242510908"
$wb.Worksheets.Item(2).Name = "NB_TO-165125552534096"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555253469613"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555254039638"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555254799607"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555242200916.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555242340908.csv"
$ws1.Range("B4").Value = "go_stims-16512555242360907.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555242500901.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_9-16512555243930907.csv"
$ws2.Range("B3").Value = "OB-16512555248879654.csv"
$ws2.Range("B4").Value = "TB-165125552532796.csv"
$ws2.Range("B5").Value = "ZB-match_0-16512555245829608.csv"
$ws2.Range("B6").Value = "TB-1651255525228961.csv"
$ws2.Range("B7").Value = "OB-16512555251699667.csv"
$ws2.Range("B8").Value = "TB-16512555252549627.csv"
$ws2.Range("B9").Value = "OB-16512555249529595.csv"
$ws2.Range("B10").Value = "ZB-match_4-16512555247449646.csv"

# --- Sheet 3: RS_TO --- (no cell content changes, only the rename above)

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555253709655.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555253499615.csv"
$ws4.Range("B4").Value = "MM_stims-16512555253869631.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255525371964.csv"
$ws4.Range("B6").Value = "MM_stims-16512555254029608.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555253879635.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512555254489622.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555254339607.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555254649634.csv"
$ws5.Range("B5").Value = "SAT_stims-1651255525407964.csv"
